$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '90.810.55'
$ws.Range("E2").Value = '  +0.60%  '
# Row 3
$ws.Range("D3").Value = '3.186.44'
$ws.Range("E3").Value = '  +2.63%  '
# Row 4
$ws.Range("E4").Value = '  +0.44%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '217.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.73%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '627.06'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.91%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.05'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +26.78%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.371'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.69%  '
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.00'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.11%  '
# Row 10
$ws.Range("D10").Value = '3.184.71'
$ws.Range("E10").Value = '  +2.64%  '
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.744'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +20.11%  '
# Row 12
$ws.Range("E12").Value = '  +6.47%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000250'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.72%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +7.60%  '
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.54'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.15%  '
# Row 16
$ws.Range("D16").Value = '90.660.49'
$ws.Range("E16").Value = '  +0.78%  '
# Row 17
$ws.Range("D17").Value = '3.781.33'
$ws.Range("E17").Value = '  +2.85%  '
# Row 18
$ws.Range("D18").Value = '3.216.85'
$ws.Range("E18").Value = '  +3.54%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +9.42%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000218'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.62%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.83%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '442.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.50%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.87%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.22'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.29%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.09'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.87%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '87.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.52%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.34'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.17%  '
# Row 28
$ws.Range("D28").Value = '3.339.11'
$ws.Range("E28").Value = '  +2.26%  '
# Row 29
$ws.Range("E29").Value = '  +0.22%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.163'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.31%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.31'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +14.10%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.984'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.59%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '528.99'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.15%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.16'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +12.02%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.74'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.45%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.62%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.145'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +11.58%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.90'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.40%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.30'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.85%  '
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.38%  '
# Row 41
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.160'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +18.41%  '
# Row 42
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '
# Row 43
$ws.Range("B43").Value = 'Hedera'
$ws.Range("C43").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0844'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +19.74%  '
# Row 44
$ws.Range("B44").Value = 'PolygonEcosystemToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.410'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.87%  '
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.05%  '
# Row 46
$ws.Range("E46").Value = '  +0.04%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '148.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.99%  '
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.34'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +9.67%  '
# Row 49
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '44.29'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.36%  '
# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '173.04'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.06%  '
# Row 51
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.38'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +9.94%  '
